$wb = $excel.ActiveWorkbook

# --- Grab sheet references ---
$ws1 = $wb.Worksheets.Item("Select Input")
$ws2 = $wb.Worksheets.Item("Radio Buttons Demo")
$ws3 = $wb.Worksheets.Item("Simple Form Demo")

# --- Add new sheet "RowColumnTable" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "RowColumnTable"
$newSheet.Cells.Item(2, 1).Value = "Sonya Frost Software Engineer Edinburgh 23 2008/12/13 `$103,600"
$newSheet.Cells.Item(1, 1).Value = "ExpectedField"
$newSheet.Range("B6").Select()

# --- Update "Select Input" sheet ---
$ws1.Range("B3").Value = "All selected colors are : Green"

# --- Set selections to match final state ---
$ws2.Range("A4").Select()
$ws3.Range("C1").Select()

# Select and activate "Select Input" sheet last since it is the active tab in the final workbook
$ws1.Activate()
$ws1.Range("B3").Select()
